# This script swaps the full contents of several row-pairs on the active
# worksheet. The underlying source data rows were re-ordered (two records
# exchanged positions) and this reproduces that by exchanging every cell
# A:AY between each pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose entire contents (columns A through AY) must be swapped.
$pairs = @(
    @(3, 4),
    @(5, 6),
    @(10, 11),
    @(14, 15),
    @(16, 17),
    @(26, 28),
    @(29, 30)
)

# Columns whose text values could be misinterpreted (auto-converted to a
# number/date) by Excel when written back via .Value2. Forcing these cells
# to the "Text" number format before assignment keeps them as plain text,
# matching the original inline-string cell types.
$textColumns = @("I", "Y", "AA")

$lastCol = "AY"

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("A" + $r1 + ":" + $lastCol + $r1)
    $rng2 = $ws.Range("A" + $r2 + ":" + $lastCol + $r2)

    # Capture current contents of both rows before overwriting either one.
    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    # Protect numeric/date-looking text cells from Excel's automatic type
    # coercion so they remain stored as text after the swap.
    foreach ($col in $textColumns) {
        $ws.Range($col + $r1).NumberFormat = "@"
        $ws.Range($col + $r2).NumberFormat = "@"
    }

    # Perform the swap.
    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}

Write-Host "Row swaps complete."
